# "full c code for motor drivers to work with pwm"
# The PWM duty-cycle related input (gear/pulley radius, column F, row 1 of
# each of the three stacked calc blocks) is lowered from 0.01 m to 0.008 m.
# Every other touched cell (F10:F15, F27:F32, F44:F49) is a formula that
# depends on it, so they recompute automatically once the input changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F1").Value = 0.008
$ws.Range("F18").Value = 0.008
$ws.Range("F35").Value = 0.008

# Restore the cursor/selection position left behind in the saved view state.
[void]$ws.Range("F49").Select()
